# Adding template with sheets
# - Rename Sheet1 -> "Invoice Screen" and Sheet2 -> "Credit Note Screen" (Sheet3 stays as-is)
# - Populate a header row on the (now) "Invoice Screen" sheet with the QA report columns
# - Leave the active selection on F1, just past the last populated header cell

$wb = $excel.ActiveWorkbook

# Rename the first two sheets; Sheet3 is left untouched.
$wb.Worksheets.Item(1).Name = "Invoice Screen"
$wb.Worksheets.Item(2).Name = "Credit Note Screen"

# Fill in the header row for the QA report template on the Invoice Screen sheet.
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A1").Value = "Test"
$ws1.Range("B1").Value = "Steps"
$ws1.Range("C1").Value = "Expected Result"
$ws1.Range("D1").Value = "Pass/Fail"
$ws1.Range("E1").Value = "Status (Fixed/Verified/Closed)"

# Move the active selection to the cell right after the populated header.
$ws1.Range("F1").Select() | Out-Null
